{"js": "// \"schedule and hist guides\"\n// Highlight the \"Histograms with JASP (made this guide)\" entry (in the\n// Guides Edits column of the schedule table) in yellow, matching the\n// highlighting already used for the other guide entries in that table.\n\nconst target = \"Histograms with JASP (made this guide)\";\n\nconst results = context.document.body.search(target, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(`Could not find text: ${target}`);\n}\n\nconst range = results.items[0];\nrange.font.highlightColor = \"Yellow\";\n\nawait context.sync();\n", "ps1": "# \"schedule and hist guides\"\n# Highlight the \"Histograms with JASP (made this guide)\" entry (in the\n# Guides Edits column of the schedule table) in yellow, matching the\n# highlighting already used for the other guide entries in that table.\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$found = $rng.Find.Execute(\"Histograms with JASP (made this guide)\")\n\nif (-not $found) {\n    throw \"Could not find target text in document\"\n}\n\n$rng.Font.HighlightColorIndex = \"Yellow\"\n"}
